$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 8: Beam 1 results ---
$ws.Range("D8").Value = 161
$ws.Range("E8").Value = 62208
$ws.Range("I8").Value = 0.000426828
$ws.Range("L8").Value = 16.415247
$ws.Range("M8").Value = -0.000000264066
$ws.Range("N8").Value = -2.374583756
$ws.Range("O8").Value = -3.82424215

# --- Row 9: Beam 2 results ---
$ws.Range("D9").Value = 135
$ws.Range("E9").Value = 62224
$ws.Range("I9").Value = 0.000426828385821
$ws.Range("J9").Value = 0.0001
$ws.Range("K9").Value = 0.0001
$ws.Range("L9").Value = 16.4152433535667
$ws.Range("M9").Value = -0.0000000024861110000000001
$ws.Range("N9").Value = -2.37458354
$ws.Range("O9").Value = -3.8242424

# --- Row 10: Beam 3 results ---
$ws.Range("D10").Value = 126
$ws.Range("E10").Value = 48768
$ws.Range("I10").Value = -0.0004519941567
$ws.Range("J10").Value = 0.00013805937
$ws.Range("K10").Value = 0.000145985128
$ws.Range("L10").Value = 17.125554
$ws.Range("M10").Value = -0.004519994
$ws.Range("N10").Value = -1.62406568
$ws.Range("O10").Value = -3.185490202

# Widen the exported columns (I:K) slightly so the new values fit
$ws.Range("I1:K1").EntireColumn.ColumnWidth = 9.28515625

# Move the active selection like the author left it
$ws.Range("J11").Select()
